$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D column values stay as text (they look numeric but must
# remain formatted exactly like the source strings, e.g. "1.00").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "48.302.83"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "2.534.84"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "324.32"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "109.38"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +4.48%  "
$ws.Range("D10").Value = "40.86"
$ws.Range("E10").Value = "  +4.08%  "
$ws.Range("D11").Value = "20.53"
$ws.Range("E11").Value = "  +11.51%  "
$ws.Range("D12").Value = "0.0832"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").Value = "7.32"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "2.926.30"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "2.539.49"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "0.862"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "48.137.73"
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("D19").Value = "13.32"
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "72.51"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("D24").Value = "271.56"
$ws.Range("E24").Value = "  +9.67%  "
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "26.34"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "0.145"
$ws.Range("E29").Value = "  +4.40%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.21"
$ws.Range("E30").Value = "  -3.80%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "35.80"
$ws.Range("E31").Value = "  +1.35%  "

$ws.Range("D32").Value = "49.78"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "19.94"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "5.43"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").Value = "4.76"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("D39").Value = "3.02"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "0.113"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "22.39"
$ws.Range("E41").Value = "  +5.56%  "
$ws.Range("D42").Value = "119.44"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").Value = "2.20"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D44").Value = "0.0301"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "2.016.34"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "3.16"
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  +4.87%  "
$ws.Range("D49").Value = "9.17"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").Value = "5.29"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").Value = "79.97"
$ws.Range("E51").Value = "  +2.54%  "
